$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new test results to row 9 (split0008): ARC_x264_Q05, ARC_x265_Q05, ARC_x265_Q10 all "ok"
$ws.Range("E9").Value = "ok"
$ws.Range("F9").Value = "ok"
$ws.Range("G9").Value = "ok"

# Add new test result to row 12 (split0011): ARC_x265_Q05 = "ok"
$ws.Range("F12").Value = "ok"

# Update the active selection to reflect the last edited cell
$ws.Activate()
$ws.Range("F12").Select()
